# "Added last minute updates" - AFMC PGI 5325
#
# 1. Paragraph 1 (the hidden "**ID__...__ID**" placeholder paragraph):
#    - add a paragraph border (top/left/bottom/right) with 5 twips of
#      space and no line (matches <w:pBdr><w:top w:space="5"/>...)
#    - widen the left indent from 120 -> 225 twips (6pt -> 11.25pt)
#    - drop the trailing " " run, leaving only the ID run
#    - rename the placeholder id from AFFARS_pgi_5325_topic_2 to
#      AFFARS_AFMC_PGI_5325
$d = $word.ActiveDocument

$para1 = $d.Paragraphs.Item(1)

# --- pBdr: space=5 on all four sides ---
$para1.Format.Borders.DistanceFromTop = 5
$para1.Format.Borders.DistanceFromLeft = 5
$para1.Format.Borders.DistanceFromBottom = 5
$para1.Format.Borders.DistanceFromRight = 5

# --- ind: left 120 -> 225 twips (twips / 20 = points) ---
$para1.Format.LeftIndent = 11.25

# --- drop the trailing " " run so only the ID run remains ---
$r1 = $para1.Range
$fullText = $r1.Text
$trimmedLen = $fullText.TrimEnd().Length
if ($trimmedLen -lt $fullText.Length) {
    $trailing = $d.Range($r1.Start + $trimmedLen, $r1.Start + $fullText.Length - 1)
    $trailing.Delete()
}

# --- rewrite the placeholder id text ---
$null = $d.Content.Find.Execute("**ID__AFFARS_pgi_5325_topic_2__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5325__ID**", 2)
